# Auto update Excel log
# Appends newly-logged sensor/alert rows to four sheets:
#   ALERTS      : rows 5-6   (A1:F4 -> A1:F6)
#   Humidity    : rows 45-54 (A1:F44 -> A1:F54)
#   Temperature : rows 45-54 (A1:F44 -> A1:F54)
#   Proximity   : rows 38-40 (A1:F37 -> A1:F40)
#
# All data in this workbook is stored as plain text (no sheet uses numeric
# cell types), so every new cell is force-formatted as Text ("@") before the
# value is written and reset back to the default "Normal" style afterwards.
# This stops Excel's autosense from turning strings like "2026-02-01" into a
# date serial or "77.8%" into a percentage number.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$RowNum,
        [string[]]$Values
    )

    $rng = $ws.Range("A$RowNum`:F$RowNum")
    $rng.NumberFormat = "@"

    $ws.Cells.Item($RowNum, 1).Value = $Values[0]
    $ws.Cells.Item($RowNum, 2).Value = $Values[1]
    $ws.Cells.Item($RowNum, 3).Value = $Values[2]
    $ws.Cells.Item($RowNum, 4).Value = $Values[3]
    $ws.Cells.Item($RowNum, 5).Value = $Values[4]
    $ws.Cells.Item($RowNum, 6).Value = $Values[5]

    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# ALERTS sheet: rows 5-6
# ---------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")

Add-LogRow $wsAlerts 5 @("2026-02-01", "18:22:36", "18:00", "Bathroom", "MINIMAL", "MINIMAL ALERT: Bathroom occupied, no motion > 20s.")
Add-LogRow $wsAlerts 6 @("2026-02-01", "18:22:57", "18:00", "Bathroom", "MODERATE", "MODERATE ALERT: Bathroom occupied, no motion > 40s.")

# ---------------------------------------------------------------------
# Humidity sheet: rows 45-54
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

Add-LogRow $wsHumidity 45 @("2026-02-01", "18:22:08", "18:00", "Bathroom", "77.8%", "Active")
Add-LogRow $wsHumidity 46 @("2026-02-01", "18:22:11", "18:00", "Bathroom", "78.9%", "Active")
Add-LogRow $wsHumidity 47 @("2026-02-01", "18:22:22", "18:00", "Bathroom", "78.8%", "Active")
Add-LogRow $wsHumidity 48 @("2026-02-01", "18:22:32", "18:00", "Bathroom", "78.8%", "Active")
Add-LogRow $wsHumidity 49 @("2026-02-01", "18:22:37", "18:00", "Bathroom", "77.9%", "Active")
Add-LogRow $wsHumidity 50 @("2026-02-01", "18:22:42", "18:00", "Bathroom", "79.1%", "Active")
Add-LogRow $wsHumidity 51 @("2026-02-01", "18:22:47", "18:00", "Bathroom", "77.8%", "Active")
Add-LogRow $wsHumidity 52 @("2026-02-01", "18:22:57", "18:00", "Bathroom", "78.0%", "Active")
Add-LogRow $wsHumidity 53 @("2026-02-01", "18:23:02", "18:00", "Bathroom", "78.7%", "Active")
Add-LogRow $wsHumidity 54 @("2026-02-01", "18:23:07", "18:00", "Bathroom", "78.0%", "Active")

# ---------------------------------------------------------------------
# Temperature sheet: rows 45-54
# ---------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")

Add-LogRow $wsTemperature 45 @("2026-02-01", "18:22:08", "18:00", "Bathroom", "29.4C", "Active")
Add-LogRow $wsTemperature 46 @("2026-02-01", "18:22:12", "18:00", "Bathroom", "29.4C", "Active")
Add-LogRow $wsTemperature 47 @("2026-02-01", "18:22:22", "18:00", "Bathroom", "29.4C", "Active")
Add-LogRow $wsTemperature 48 @("2026-02-01", "18:22:32", "18:00", "Bathroom", "29.4C", "Active")
Add-LogRow $wsTemperature 49 @("2026-02-01", "18:22:37", "18:00", "Bathroom", "29.4C", "Active")
Add-LogRow $wsTemperature 50 @("2026-02-01", "18:22:42", "18:00", "Bathroom", "29.4C", "Active")
Add-LogRow $wsTemperature 51 @("2026-02-01", "18:22:47", "18:00", "Bathroom", "29.4C", "Active")
Add-LogRow $wsTemperature 52 @("2026-02-01", "18:22:57", "18:00", "Bathroom", "29.4C", "Active")
Add-LogRow $wsTemperature 53 @("2026-02-01", "18:23:02", "18:00", "Bathroom", "29.3C", "Active")
Add-LogRow $wsTemperature 54 @("2026-02-01", "18:23:07", "18:00", "Bathroom", "29.4C", "Active")

# ---------------------------------------------------------------------
# Proximity sheet: rows 38-40
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")

Add-LogRow $wsProximity 38 @("2026-02-01", "18:22:07", "18:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom")
Add-LogRow $wsProximity 39 @("2026-02-01", "18:22:08", "18:00", "Bathroom Door", "EXIT", "User EXITED Bathroom")
Add-LogRow $wsProximity 40 @("2026-02-01", "18:22:14", "18:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom")
